# Update odds/stats values on Sheet1 to match the latest FlashScore export.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.5
$ws.Range("I2").Value = 3
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("Q2").Value = 2.6
$ws.Range("R2").Value = 1.48
$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 23
$ws.Range("AK2").Value = 34
$ws.Range("AN2").Value = 4.33
$ws.Range("AV2").Value = 81

# Row 4
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5

# Row 5
$ws.Range("G5").Value = 1.25
$ws.Range("H5").Value = 5
$ws.Range("K5").Value = 2.3
$ws.Range("L5").Value = 13
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25
$ws.Range("Q5").Value = 2.08
$ws.Range("R5").Value = 1.73
$ws.Range("Y5").Value = 11
$ws.Range("Z5").Value = 6.5
$ws.Range("AB5").Value = 51
$ws.Range("AH5").Value = 23
$ws.Range("AL5").Value = 151
$ws.Range("AS5").Value = 351

# Row 6
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10

# Row 8
$ws.Range("Q8").Value = 1.9
$ws.Range("R8").Value = 1.95

# Row 9
$ws.Range("Q9").Value = 2.05
$ws.Range("R9").Value = 1.75

# Row 10
$ws.Range("G10").Value = 1.65
$ws.Range("I10").Value = 5.75
$ws.Range("J10").Value = 2.3
$ws.Range("K10").Value = 2.05
$ws.Range("L10").Value = 6
$ws.Range("N10").Value = 7.5
$ws.Range("W10").Value = 5.5
$ws.Range("AC10").Value = 7.5
$ws.Range("AI10").Value = 29
$ws.Range("AK10").Value = 67
$ws.Range("AL10").Value = 51
$ws.Range("AW10").Value = 7
$ws.Range("AX10").Value = 34
